$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update frequency values in column C (rows 2-10)
$ws.Range("C2").Value = 2980
$ws.Range("C3").Value = 2943
$ws.Range("C4").Value = 2174
$ws.Range("C5").Value = 1327
$ws.Range("C6").Value = 1171
$ws.Range("C7").Value = 663
$ws.Range("C8").Value = 584
$ws.Range("C9").Value = 434
$ws.Range("C10").Value = 428

# Row 11: swap the category values between A11 and B11
$ws.Range("A11").Value = "Home Decor"
$ws.Range("B11").Value = "Seasonal & Holidays"
$ws.Range("C11").Value = 401
